$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCT1")
$ws.Activate()

# "Private cars_Maximum stock" (row 7) is derived from the base-year
# Private Car stock (row 42) with a +5% factor; bump the factor to +6%.
$ws.Range("H7").Formula = "=H42*1.06"
$ws.Range("I7:AH7").Formula = "=I42*1.06"

# Update the sheet selection / view to match the author's saved state.
$ws.Range("M20").Select()
